$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 11, shifting existing rows 11..109 down to 12..110
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new data record
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(11, 3).Value = "Metropolitana"
$ws.Cells.Item(11, 4).Value = 44490
$ws.Cells.Item(11, 5).Value = 13
$ws.Cells.Item(11, 6).Value = "Fruta"
$ws.Cells.Item(11, 7).Value = 100101
$ws.Cells.Item(11, 8).Value = "Berries"
$ws.Cells.Item(11, 9).Value = 100101001
$ws.Cells.Item(11, 10).Value = "Arándano (blue)"
$ws.Cells.Item(11, 11).Value = "Sin especificar"
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 330
$ws.Cells.Item(11, 14).Value = 12000
$ws.Cells.Item(11, 15).Value = 12000
$ws.Cells.Item(11, 16).Value = 12000
$ws.Cells.Item(11, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(11, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(11, 19).Value = 6000
$ws.Cells.Item(11, 20).Value = 2
